# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) to the results sheet and updates
# recomputed metric values (columns D/E/F) from the refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header - copy the header formatting (bold, border, centered)
# from the neighboring header cell, then set the text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# --- Block 1 (Iterations = 100), rows 2-11 ---

# Row 2 - Control 26
$ws.Range("D2").Value = 0.4807142772271398
$ws.Range("E2").Value = 0.4807142772271398
$ws.Range("H2").Value = 0

# Row 3 - Control 33
$ws.Range("D3").Value = 0.6235864263168943
$ws.Range("E3").Value = 0.6235864263168943
$ws.Range("H3").Value = 0

# Row 4 - Control 36
$ws.Range("D4").Value = 0.6298392979282658
$ws.Range("E4").Value = 0.6298392979282658
$ws.Range("H4").Value = 0

# Row 5 - Control 49
$ws.Range("D5").Value = 0.2530959129238339
$ws.Range("E5").Value = 0.2530959129238339
$ws.Range("H5").Value = 0

# Row 6 - Control 2
$ws.Range("D6").Value = 0.62511108607642
$ws.Range("E6").Value = 0.62511108607642
$ws.Range("H6").Value = 0

# Row 7 - MDD 0
$ws.Range("D7").Value = 0.6302063365315793
$ws.Range("E7").Value = 0.3697936634684207
$ws.Range("H7").Value = 1

# Row 8 - MDD 30
$ws.Range("D8").Value = 0.6121140773674206
$ws.Range("E8").Value = 0.3878859226325794
$ws.Range("H8").Value = 1

# Row 9 - MDD 46
$ws.Range("D9").Value = 0.622641324101204
$ws.Range("E9").Value = 0.377358675898796
$ws.Range("H9").Value = 1

# Row 10 - MDD 17
$ws.Range("D10").Value = 0.6094490706054323
$ws.Range("E10").Value = 0.3905509293945677
$ws.Range("H10").Value = 1

# Row 11 - MDD 23
$ws.Range("D11").Value = 0.6206902391659949
$ws.Range("E11").Value = 0.3793097608340051
$ws.Range("F11").Value = 0.6297582387924194
$ws.Range("H11").Value = 1

# --- Block 2 (Iterations = 200), rows 12-21 ---
# D/E/F values unchanged for this block; only the new Label column is added.

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
